$d = $word.ActiveDocument

# 1) Fill in the previously-empty second paragraph with the GitHub link.
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Text = "https://github.com/brayden276/Brayden-SIT323"

# 2) Mark the run containing the screenshot drawing as NoProof.
$p4 = $d.Paragraphs.Item(4)
$p4.Range.NoProofing = $true
